$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.645.49"
$ws.Range("E2").Value = "'  +0.82%  "
$ws.Range("D3").Value = "'1.636.62"
$ws.Range("E3").Value = "'  -0.36%  "
$ws.Range("E4").Value = "'  +0.13%  "
$ws.Range("D5").Value = "'212.64"
$ws.Range("E5").Value = "'  +0.26%  "
$ws.Range("E6").Value = "'  -1.07%  "
$ws.Range("E7").Value = "'  +0.14%  "
$ws.Range("D8").Value = "'22.96"
$ws.Range("E8").Value = "'  -0.58%  "
$ws.Range("E9").Value = "'  +0.48%  "
$ws.Range("E10").Value = "'  -0.12%  "
$ws.Range("D11").Value = "'0.0894"
$ws.Range("E11").Value = "'  +0.66%  "
$ws.Range("D12").Value = "'1.870.57"
$ws.Range("D13").Value = "'1.639.33"
$ws.Range("E13").Value = "'  -0.12%  "
$ws.Range("E14").Value = "'  -0.04%  "
$ws.Range("D15").Value = "'0.560"
$ws.Range("E15").Value = "'  -1.99%  "
$ws.Range("D16").Value = "'64.57"
$ws.Range("E16").Value = "'  +0.20%  "
$ws.Range("D17").Value = "'27.646.20"
$ws.Range("E17").Value = "'  +0.89%  "
$ws.Range("D18").Value = "'229.67"
$ws.Range("E18").Value = "'  +0.03%  "
$ws.Range("D19").Value = "'7.75"
$ws.Range("E19").Value = "'  +1.87%  "
$ws.Range("E20").Value = "'  +0.11%  "
$ws.Range("E21").Value = "'  +0.18%  "
$ws.Range("E22").Value = "'  -1.13%  "
$ws.Range("E23").Value = "'  +3.92%  "
$ws.Range("E24").Value = "'  -2.96%  "
$ws.Range("D25").Value = "'150.19"
$ws.Range("E25").Value = "'  +2.09%  "
$ws.Range("E26").Value = "'  -1.16%  "
$ws.Range("E27").Value = "'  -1.52%  "
$ws.Range("D28").Value = "'15.64"
$ws.Range("E28").Value = "'  +0.41%  "
$ws.Range("E29").Value = "'  +0.03%  "
$ws.Range("D30").Value = "'1.18"
$ws.Range("E30").Value = "'  +0.36%  "
$ws.Range("E31").Value = "'  +0.07%  "
$ws.Range("E32").Value = "'  +0.22%  "
$ws.Range("D33").Value = "'1.453.17"
$ws.Range("E33").Value = "'  +2.85%  "
$ws.Range("E34").Value = "'  -1.46%  "
$ws.Range("E35").Value = "'  -0.94%  "
$ws.Range("E36").Value = "'  +0.41%  "
$ws.Range("E37").Value = "'  -0.03%  "
$ws.Range("D38").Value = "'0.875"
$ws.Range("E38").Value = "'  -1.20%  "
$ws.Range("D39").Value = "'0.0166"
$ws.Range("E39").Value = "'  +0.14%  "
$ws.Range("D40").Value = "'0.899"
$ws.Range("E40").Value = "'  +9.57%  "
$ws.Range("D41").Value = "'69.75"
$ws.Range("E41").Value = "'  +7.99%  "
$ws.Range("E42").Value = "'  -0.80%  "
$ws.Range("E43").Value = "'  +0.16%  "
$ws.Range("E44").Value = "'  +1.61%  "
$ws.Range("E45").Value = "'  +0.35%  "
$ws.Range("E46").Value = "'  -0.24%  "
$ws.Range("D47").Value = "'1.779.77"
$ws.Range("E47").Value = "'  -0.24%  "
$ws.Range("E48").Value = "'  +2.19%  "
$ws.Range("D49").Value = "'86.13"
$ws.Range("E49").Value = "'  -2.12%  "
$ws.Range("D50").Value = "'0.0₆0106"
$ws.Range("E50").Value = "'  -1.20%  "
$ws.Range("D51").Value = "'0.0985"
$ws.Range("E51").Value = "'  -0.59%  "
